$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.690.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.828.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4662'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3607'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07145'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9049'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07776'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.44'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.831.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.272'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.352'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.62%  '
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008575'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.741.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.016'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.924'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.974'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.833'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08810'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.149'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7338'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.149'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.99%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.739'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.450'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.079'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01925'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.932'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05133'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.882'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5071'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1499'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.059'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('E44').Value = '  +1.00%  '
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('E46').Value = '  +2.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.37'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.565'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06058'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('E51').Value = '  -0.33%  '
